# Auto commit at 2025-12-22 8:56:02.27
# Updates the "Metrics" sheet's monthly-to-date figures (column B, rows 2-13)
# with refreshed totals. The "today" sheet pulls these via formulas
# (Metrics!B2 ... Metrics!B13) and will recalculate automatically.

$wb = $excel.ActiveWorkbook

# --- Metrics sheet: refresh column B values -------------------------------
$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value  = 280952.52
$metrics.Range("B3").Value  = 240372.52
$metrics.Range("B4").Value  = 86055.9
$metrics.Range("B5").Value  = 11509
$metrics.Range("B6").Value  = 5483659.6300000008
$metrics.Range("B7").Value  = 4640725.4800000004
$metrics.Range("B8").Value  = 1618012.7800000003
$metrics.Range("B9").Value  = 214216
$metrics.Range("B10").Value = 33949040.61999999
$metrics.Range("B11").Value = 31916000.639999997
$metrics.Range("B12").Value = 11899734.819999995
$metrics.Range("B13").Value = 1311846

# Restore the selection that was active on this sheet.
[void]$metrics.Range("D10").Select()

# --- "today" sheet: just move the selection --------------------------------
# (B11:B22/E/F columns are formulas referencing Metrics!B2:B13 and update
# automatically on recalculation; A1 is the volatile TODAY()-1 cell and is
# left as-is so it keeps recalculating from the live clock.)
$today = $wb.Worksheets.Item("today")
[void]$today.Range("D8").Select()
